$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet: Debug_wbSaveReport flips from FALSE to TRUE
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("B9").Value = $true

# ---------------------------------------------------------------------
# Workblocks sheet: add a new "RecoverApps" workblock pair and rename
# the existing workblock values to short labels
# ---------------------------------------------------------------------
$wsWorkblocks = $wb.Worksheets.Item("Workblocks")

# Insert two new rows right after the wbInit_SuppressSuccessful row (row 4)
$wsWorkblocks.Rows.Item(5).Insert()
$wsWorkblocks.Rows.Item(5).Insert()

# Copy formatting from the wbInit_Type / wbInit_SuppressSuccessful rows so
# the new rows get the same cell styles (font/alignment)
$wsWorkblocks.Range("A3:C4").Copy()
$wsWorkblocks.Range("A5:C6").PasteSpecial(-4122)
$wsWorkblocks.Application.CutCopyMode = $false

$wsWorkblocks.Range("A5").Value = "wbCloseAppsRecover_Type"
$wsWorkblocks.Range("B5").Value = "RecoverApps"
$wsWorkblocks.Range("C5").Value = "Name of Workblock"

$wsWorkblocks.Range("A6").Value = "wbCloseAppsRecover_SuppressSuccessful"
$wsWorkblocks.Range("B6").Value = $true
$wsWorkblocks.Range("C6").Value = "Do not log successful executions of wb"

# Shorten the existing workblock type values
$wsWorkblocks.Range("B3").Value = "Init"
$wsWorkblocks.Range("B7").Value = "GetData"
$wsWorkblocks.Range("B9").Value = "Process"
$wsWorkblocks.Range("B11").Value = "Next"
$wsWorkblocks.Range("B13").Value = "CloseApps"
$wsWorkblocks.Range("B15").Value = "InitApps"
$wsWorkblocks.Range("B17").Value = "ProcessApps"

# Widen column B to fit the new header/content layout
$wsWorkblocks.Columns.Item(2).ColumnWidth = 32.67

# ---------------------------------------------------------------------
# Tasks sheet: rename task labels and enable the FirstRun task
# ---------------------------------------------------------------------
$wsTasks = $wb.Worksheets.Item("Tasks")
$wsTasks.Range("B3").Value = "FirstRunTask"
$wsTasks.Range("B5").Value = $true
$wsTasks.Range("B6").Value = "GetDataTask"

$wsTasks.Activate()
$wsTasks.Range("B5").Select()
